$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.874799999999993
$ws.Range("D4").Value = -7.126799999999995
$ws.Range("B6").Value = 4.747200000000004
$ws.Range("B7").Value = 5.512700000000001
$ws.Range("D9").Value = -7.562199999999995
$ws.Range("D12").Value = -5.883899999999995
$ws.Range("B16").Value = 4.791599999999999
$ws.Range("D17").Value = -8.04409999999999
$ws.Range("D18").Value = -9.230199999999989
$ws.Range("D19").Value = -8.56819999999999
$ws.Range("B20").Value = 9.102199999999987
$ws.Range("D20").Value = -7.989599999999988
$ws.Range("D26").Value = -7.785300000000007
$ws.Range("B28").Value = 5.967900000000006
$ws.Range("B29").Value = 5.3195
$ws.Range("D31").Value = -7.669899999999995
$ws.Range("B32").Value = 6.615
$ws.Range("D39").Value = -8.239199999999999
$ws.Range("B40").Value = 9.587399999999995
$ws.Range("D40").Value = -8.755399999999991
$ws.Range("D41").Value = -7.64129999999999
$ws.Range("D42").Value = -8.286699999999989
$ws.Range("D43").Value = -7.468500000000003
$ws.Range("B46").Value = 5.838100000000002
$ws.Range("D47").Value = -7.628600000000002
$ws.Range("D48").Value = -7.424699999999996
$ws.Range("B51").Value = 5.987300000000001
$ws.Range("B52").Value = 5.627599999999997
$ws.Range("B57").Value = 4.959599999999997
$ws.Range("B59").Value = 4.931000000000001
$ws.Range("B62").Value = 5.406199999999996
$ws.Range("D63").Value = -6.592899999999994
$ws.Range("D64").Value = -6.960699999999991
$ws.Range("B66").Value = 5.678099999999993
$ws.Range("B73").Value = 9.082000000000004
$ws.Range("B74").Value = 9.055799999999994
$ws.Range("D76").Value = -7.042299999999996
$ws.Range("D81").Value = -7.721299999999998
$ws.Range("D89").Value = -8.290999999999995
$ws.Range("B92").Value = 5.052599999999995
$ws.Range("D94").Value = -5.923599999999998
$ws.Range("B100").Value = 4.797800000000002
